# Metadata.xlsx test-fixture update
# - Volume_Unit "L" -> "µL" for every data row
# - Norm1_Unit "cell" -> "million_cells" for every data row
# - Resuspension_Volume values 2.5E-5 -> 200, and cleared to the
#   default (unstyled / General) cell style (was scientific-notation)
# - Norm1 values 1000000 -> 0.6
# - Reselect E2:E23 (new active selection) instead of H12
# - Autofit column E now that its content changed width

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 23

# --- Column C: Volume_Unit ------------------------------------------------
$ws.Range("C2:C$lastRow").Value = "µL"

# --- Column E: Norm1_Unit ---------------------------------------------------
$ws.Range("E2:E$lastRow").Value = "million_cells"

# --- Column B: Resuspension_Volume ------------------------------------------
$ws.Range("B2:B$lastRow").Value = 200
# Drop the scientific-notation number format these cells had (numFmtId 11)
# back to the workbook's default / unstyled look.
$ws.Range("B2:B$lastRow").Style = "Normal"

# --- Column D: Norm1 ----------------------------------------------------
$ws.Range("D2:D$lastRow").Value = 0.6

# --- Column D no longer needs an explicit best-fit width; column E does,
#     since "million_cells" is wider than "cell".
$ws.Columns.Item(5).AutoFit()

# --- Selection moves from the old H12 cell to the updated E2:E23 range -----
$ws.Range("E2:E$lastRow").Select()
